$d = $word.ActiveDocument

$d.Content.Find.Execute("921÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "921÷7=", 2) | Out-Null
$d.Content.Find.Execute("701÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "663÷4=", 2) | Out-Null
$d.Content.Find.Execute("975÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "726÷8=", 2) | Out-Null
$d.Content.Find.Execute("316÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "664÷6=", 2) | Out-Null
$d.Content.Find.Execute("234÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "404÷9=", 2) | Out-Null
$d.Content.Find.Execute("355÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "322÷6=", 2) | Out-Null
$d.Content.Find.Execute("143÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "131÷5=", 2) | Out-Null
$d.Content.Find.Execute("493÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "497÷4=", 2) | Out-Null
$d.Content.Find.Execute("252÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "982÷2=", 2) | Out-Null
$d.Content.Find.Execute("912÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "317÷3=", 2) | Out-Null
$d.Content.Find.Execute("970÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "634÷9=", 2) | Out-Null
$d.Content.Find.Execute("684÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "849÷4=", 2) | Out-Null
$d.Content.Find.Execute("893÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "206÷3=", 2) | Out-Null
$d.Content.Find.Execute("687÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "579÷5=", 2) | Out-Null
$d.Content.Find.Execute("406÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "108÷6=", 2) | Out-Null
$d.Content.Find.Execute("515÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "825÷9=", 2) | Out-Null
$d.Content.Find.Execute("704÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "980÷5=", 2) | Out-Null
$d.Content.Find.Execute("907÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "340÷4=", 2) | Out-Null
$d.Content.Find.Execute("892÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "105÷8=", 2) | Out-Null
$d.Content.Find.Execute("254÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "435÷4=", 2) | Out-Null
$d.Content.Find.Execute("900÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "225÷3=", 2) | Out-Null
$d.Content.Find.Execute("278÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "795÷7=", 2) | Out-Null
$d.Content.Find.Execute("563÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "680÷9=", 2) | Out-Null
$d.Content.Find.Execute("978÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "940÷8=", 2) | Out-Null
$d.Content.Find.Execute("286÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "146÷2=", 2) | Out-Null
